# Added WT testcases with element util modification
#
# Target sheet: "TestData_UAT" (the WebTours / "WT" test-data table).
# Row 2/3 gain three new columns (firstname/lastname/password headers +
# sumit/raut/pwd sample values), and a brand-new data row is effectively
# cloned from row 3 into row 5 (same Execute/email/gender/firstname/
# lastname/password values, same "Hyperlink" formatting + mailto link on
# column D), replacing the old rautsumit2/rautsumit placeholder row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData_UAT")

# --- Row 2 (sub-header row): add the three new column headers ---------
$ws.Range("F2").Value = "firstname"
$ws.Range("G2").Value = "lastname"
$ws.Range("H2").Value = "password"

# --- Row 3 (first WT data row): add the matching sample values --------
$ws.Range("F3").Value = "sumit"
$ws.Range("G3").Value = "raut"
$ws.Range("H3").Value = "pwd"

# --- Row 5 (second WT data row): replace the old placeholder values ---
# First, clone D3's "Hyperlink" cell formatting onto D5 (style only;
# the value/hyperlink relationship are set explicitly below) so D5 ends
# up visually identical to D3, matching the other WT row.
$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C5").Value = "yes"
$ws.Range("D5").Value = "rautsumit@test.com"
$ws.Range("E5").Value = "male"
$ws.Range("F5").Value = "sumit"
$ws.Range("G5").Value = "raut"
$ws.Range("H5").Value = "pwd"

# Row 5 now mirrors row 3's taller height.
$ws.Rows(5).RowHeight = 45

# D5 becomes a mailto hyperlink, same as D3.
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:rautsumit@test.com")

# Selection moves down onto the newly-edited row.
$ws.Range("A5:XFD5").Select()
